# Add a new "lifeform_wright" column (AL) to the worksheet, populating the
# header in AL1 and the per-species lifeform values in AL2:AL15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - match formatting of the other header cells (bold, centered)
$ws.Range("AL1").Value = "lifeform_wright"
$ws.Range("AL1").Font.Bold = $true
$ws.Range("AL1").HorizontalAlignment = -4108  # xlCenter

$values = @{
    2  = "LIANA, CLIMBER"
    3  = "UNDERSTORY, FREE"
    4  = "SHRUB, FREE"
    5  = "SHRUB, FREE"
    6  = "SHRUB, FREE"
    7  = "SHRUB, FREE"
    8  = "SHRUB, FREE"
    9  = "SHRUB, FREE"
    10 = "SHRUB, FREE"
    11 = "SHRUB, FREE"
    12 = "UNDERSTORY, FREE"
    13 = "LIANA, CLIMBER"
    14 = "TREE, FREE"
    15 = "SHRUB, FREE"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 38).Value = $values[$row]
}
